$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("N91").ClearContents()
$ws.Range("H109").Value = 40000
$ws.Range("J109").Value = 40000
$ws.Range("L109").Value = 40000
$ws.Range("N109").Value = -42774
$ws.Range("H138").Value = 4353.9116
$ws.Range("J138").Value = 5136
$ws.Range("L138").Value = 15408
$ws.Range("N138").Value = -25688

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 167883
$ws.Range("I45").Value = 167883
$ws.Range("K45").Value = 167883
$ws.Range("M45").Value = -167506
$ws.Range("H74").Value = 751.1429000000001
$ws.Range("I74").Value = 751.1429000000001
$ws.Range("K74").Value = 751.1429000000001
$ws.Range("M74").Value = 122.8570999999999
$ws.Range("H77").Value = 751.1429000000001
$ws.Range("I77").Value = 751.1429000000001
$ws.Range("K77").Value = 3755.7145
$ws.Range("M77").Value = 612.2855

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 3000
$ws.Range("I29").Value = 3000
$ws.Range("K29").Value = 3000
$ws.Range("M29").Value = -2711
$ws.Range("H36").Value = 3993
$ws.Range("J36").Value = 4989.5
$ws.Range("L36").Value = 4989.5
$ws.Range("N36").Value = -6057.5
$ws.Range("H64").Value = 946.875
$ws.Range("J64").Value = 948
$ws.Range("L64").Value = 948
$ws.Range("N64").Value = -1398
$ws.Range("H67").Value = 946.875
$ws.Range("J67").Value = 948
$ws.Range("L67").Value = 948
$ws.Range("N67").Value = -2508
$ws.Range("H75").Value = 61666.668
$ws.Range("I75").Value = 17500
$ws.Range("K75").Value = 17500
$ws.Range("M75").Value = -16564
$ws.Range("H78").Value = 61666.668
$ws.Range("I78").Value = 17500
$ws.Range("K78").Value = 52500
$ws.Range("M78").Value = -47820
$ws.Range("H94").Value = 549.25
$ws.Range("I94").Value = 399.5
$ws.Range("K94").Value = 399.5
$ws.Range("M94").Value = 51.5
$ws.Range("H105").Value = 2586.3076
$ws.Range("I105").Value = 2446.2856
$ws.Range("J105").Value = 2749.6667
$ws.Range("K105").Value = 2446.2856
$ws.Range("L105").Value = 2749.6667
$ws.Range("M105").Value = -699.2856000000002
$ws.Range("N105").Value = -6243.6667

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 10817.667
$ws.Range("I86").Value = 11224.75
$ws.Range("J86").Value = 10003.5
$ws.Range("K86").Value = 11224.75
$ws.Range("L86").Value = 10003.5
$ws.Range("M86").Value = -10101.75
$ws.Range("N86").Value = -12249.5
$ws.Range("H89").Value = 10817.667
$ws.Range("I89").Value = 11224.75
$ws.Range("J89").Value = 10003.5
$ws.Range("K89").Value = 56123.75
$ws.Range("L89").Value = 50017.5
$ws.Range("M89").Value = -50507.75
$ws.Range("N89").Value = -61249.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 997.5
$ws.Range("I22").Value = 997.5
$ws.Range("K22").Value = 2992.5
$ws.Range("M22").Value = -2823.5
$ws.Range("H27").Value = 997.5
$ws.Range("I27").Value = 997.5
$ws.Range("K27").Value = 2992.5
$ws.Range("M27").Value = -2890.5
$ws.Range("H44").Value = 452.1
$ws.Range("I44").Value = 546.7143
$ws.Range("J44").Value = 401.15384
$ws.Range("K44").Value = 1640.1429
$ws.Range("L44").Value = 1203.46152
$ws.Range("M44").Value = -1242.1429
$ws.Range("N44").Value = -1999.46152
$ws.Range("H55").Value = 750
$ws.Range("I55").Value = 500
$ws.Range("J55").Value = 1000
$ws.Range("K55").Value = 1500
$ws.Range("L55").Value = 3000
$ws.Range("M55").Value = -1323
$ws.Range("N55").Value = -3354
$ws.Range("H113").Value = 1613.9
$ws.Range("J113").Value = 1536.75
$ws.Range("L113").Value = 4610.25
$ws.Range("N113").Value = -8950.25
$ws.Range("H132").Value = 673.6667
$ws.Range("I132").Value = 648
$ws.Range("J132").Value = 699.3333
$ws.Range("K132").Value = 5832
$ws.Range("L132").Value = 6293.9997
$ws.Range("M132").Value = -3302
$ws.Range("N132").Value = -11353.9997

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6699
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 6699
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H80").Value = 3354.6428
$ws.Range("I80").Value = 3221.375
$ws.Range("J80").Value = 3532.3333
$ws.Range("K80").Value = 3221.375
$ws.Range("L80").Value = 3532.3333
$ws.Range("M80").Value = -2223.375
$ws.Range("N80").Value = -5528.3333
$ws.Range("H83").Value = 3354.6428
$ws.Range("I83").Value = 3221.375
$ws.Range("J83").Value = 3532.3333
$ws.Range("K83").Value = 16106.875
$ws.Range("L83").Value = 17661.6665
$ws.Range("M83").Value = -11114.875
$ws.Range("N83").Value = -27645.6665
$ws.Range("H140").Value = 142997
$ws.Range("J140").Value = 142997
$ws.Range("L140").Value = 142997
$ws.Range("N140").Value = -153357
$ws.Range("H141").Value = 59999
$ws.Range("J141").Value = 59999
$ws.Range("L141").Value = 59999
$ws.Range("N141").Value = -70359

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3299
$ws.Range("I40").Value = 2565.3333
$ws.Range("K40").Value = 2565.3333
$ws.Range("M40").Value = -2429.3333
$ws.Range("H100").Value = 4102.8887
$ws.Range("I100").Value = 3349.2
$ws.Range("K100").Value = 3349.2
$ws.Range("M100").Value = -2808.2
$ws.Range("H122").Value = 5515.724
$ws.Range("I122").Value = 4629.8423
$ws.Range("J122").Value = 7198.9
$ws.Range("K122").Value = 13889.5269
$ws.Range("L122").Value = 21596.7
$ws.Range("M122").Value = -11439.5269
$ws.Range("N122").Value = -26496.7
$ws.Range("H136").Value = 2604.3125
$ws.Range("I136").Value = 1151.3636
$ws.Range("K136").Value = 3454.0908
$ws.Range("M136").Value = -904.0907999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 39999
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H81").Value = 671.6
$ws.Range("I81").Value = 539.5
$ws.Range("K81").Value = 1079
$ws.Range("M81").Value = -18
$ws.Range("H84").Value = 671.6
$ws.Range("I84").Value = 539.5
$ws.Range("K84").Value = 5395
$ws.Range("M84").Value = -91
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H141").Value = 182165.5
$ws.Range("J141").Value = 175798.8
$ws.Range("L141").Value = 175798.8
$ws.Range("N141").Value = -186158.8
